$d = $word.ActiveDocument

# --- Step 1: remove the bulk of the body content -----------------------
# Everything from "Problem:" (paragraph 2) through "Plot the prediction"
# (paragraph 20) is deleted outright, collapsing those paragraph marks too.
$midStart = $d.Paragraphs.Item(2).Range.Start
$midEnd   = $d.Paragraphs.Item(20).Range.End
$midRange = $d.Range($midStart, $midEnd)
$midRange.Delete()

# --- Step 2: collapse the title paragraph's two runs into one ----------
# "Name of" + " App" -> "daniel"
$titlePara  = $d.Paragraphs.Item(1)
$titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$titleRange.Text = "daniel"

# --- Step 3: strip the "Hard" text but keep its paragraph/bookmark -----
$hardPara  = $d.Paragraphs.Item(2)
$hardRange = $d.Range($hardPara.Range.Start, $hardPara.Range.End - 1)
$hardRange.Delete()

# --- Step 4: remove the trailing "Learn the pattern and trends" para ---
$lastPara  = $d.Paragraphs.Item(3)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Delete()
